# Resequenced columns on prop bets for better output
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: column headers now map to different shared-string
# indices after the resequencing (B=mp, C=foulvictory, D=autorun,
# E=quests, F=climbs, G=boss, H=puavg, I=powerUp)
$ws.Range("B1").Value = "mp"
$ws.Range("C1").Value = "foulvictory"
$ws.Range("D1").Value = "autorun"
$ws.Range("E1").Value = "quests"
$ws.Range("F1").Value = "climbs"
$ws.Range("G1").Value = "boss"
$ws.Range("H1").Value = "puavg"
$ws.Range("I1").Value = "powerUp"

# Update data rows with resequenced / recalculated stat values
$ws.Range("B2").Value = 80
$ws.Range("C2").Value = 0.05
$ws.Range("D2").Value = 4.25
$ws.Range("E2").Value = 0.35
$ws.Range("F2").Value = 0.3875
$ws.Range("G2").Value = 0.0125
$ws.Range("H2").Value = 2.7625
$ws.Range("I2").Value = '{''Levitate'': 102, ''Force'': [50, 18, 20, 12], ''Boost'': [69, 31, 18, 20]}'
$ws.Range("B8").Value = 80
$ws.Range("C8").Value = 0.0375
$ws.Range("D8").Value = 4.625
$ws.Range("E8").Value = 0.425
$ws.Range("F8").Value = 0.8
$ws.Range("G8").Value = 0.0375
$ws.Range("H8").Value = 2.7875
$ws.Range("I8").Value = '{''Levitate'': 105, ''Force'': [38, 14, 15, 9], ''Boost'': [80, 23, 19, 38]}'
$ws.Range("B9").Value = 35
$ws.Range("C9").Value = 0.0571428571428571
$ws.Range("D9").Value = 3.25714285714286
$ws.Range("E9").Value = 0.114285714285714
$ws.Range("F9").Value = 0.257142857142857
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 2.4
$ws.Range("I9").Value = '{''Levitate'': 20, ''Force'': [22, 6, 11, 5], ''Boost'': [42, 30, 7, 5]}'
$ws.Range("B10").Value = 87
$ws.Range("C10").Value = 0.0689655172413793
$ws.Range("D10").Value = 5.06896551724138
$ws.Range("E10").Value = 0.850574712643678
$ws.Range("F10").Value = 1.27586206896552
$ws.Range("G10").Value = 0.218390804597701
$ws.Range("H10").Value = 3.16091954022989
$ws.Range("I10").Value = '{''Levitate'': 119, ''Force'': [63, 24, 20, 19], ''Boost'': [93, 40, 22, 31]}'
$ws.Range("B12").Value = 70
$ws.Range("C12").Value = 0.0428571428571429
$ws.Range("D12").Value = 5.44285714285714
$ws.Range("E12").Value = 1.2
$ws.Range("F12").Value = 1.35714285714286
$ws.Range("G12").Value = 0.0857142857142857
$ws.Range("H12").Value = 3.51428571428571
$ws.Range("I12").Value = '{''Levitate'': 121, ''Force'': [60, 17, 19, 24], ''Boost'': [65, 16, 17, 32]}'
$ws.Range("B13").Value = 92
$ws.Range("C13").Value = 0.0434782608695652
$ws.Range("D13").Value = 5.35869565217391
$ws.Range("E13").Value = 0.967391304347826
$ws.Range("F13").Value = 1.30434782608696
$ws.Range("G13").Value = 0.0978260869565217
$ws.Range("H13").Value = 3.29347826086957
$ws.Range("I13").Value = '{''Levitate'': 138, ''Force'': [64, 20, 30, 14], ''Boost'': [101, 33, 38, 30]}'
$ws.Range("B14").Value = 56
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 4.98214285714286
$ws.Range("E14").Value = 0.589285714285714
$ws.Range("F14").Value = 0.928571428571429
$ws.Range("G14").Value = 0.0892857142857143
$ws.Range("H14").Value = 3.375
$ws.Range("I14").Value = '{''Levitate'': 84, ''Force'': [43, 14, 21, 8], ''Boost'': [62, 22, 17, 23]}'
$ws.Range("B16").Value = 74
$ws.Range("C16").Value = 0.0540540540540541
$ws.Range("D16").Value = 4.58108108108108
$ws.Range("E16").Value = 0.689189189189189
$ws.Range("F16").Value = 0.689189189189189
$ws.Range("G16").Value = 0.0405405405405405
$ws.Range("H16").Value = 2.35135135135135
$ws.Range("I16").Value = '{''Levitate'': 84, ''Force'': [26, 10, 10, 6], ''Boost'': [64, 30, 10, 24]}'
$ws.Range("B18").Value = 83
$ws.Range("C18").Value = 0.0843373493975904
$ws.Range("D18").Value = 4.79518072289157
$ws.Range("E18").Value = 0.578313253012048
$ws.Range("F18").Value = 1.32530120481928
$ws.Range("G18").Value = 0.156626506024096
$ws.Range("H18").Value = 2.90361445783133
$ws.Range("I18").Value = '{''Levitate'': 117, ''Force'': [46, 15, 22, 9], ''Boost'': [78, 33, 10, 35]}'
$ws.Range("B19").Value = 61
$ws.Range("C19").Value = 0.0163934426229508
$ws.Range("D19").Value = 4.72131147540984
$ws.Range("E19").Value = 0.524590163934426
$ws.Range("F19").Value = 0.934426229508197
$ws.Range("G19").Value = 0.0491803278688525
$ws.Range("H19").Value = 2.18032786885246
$ws.Range("I19").Value = '{''Levitate'': 62, ''Force'': [22, 9, 6, 7], ''Boost'': [49, 23, 13, 13]}'
$ws.Range("B21").Value = 78
$ws.Range("C21").Value = 0.0256410256410256
$ws.Range("D21").Value = 5.12820512820513
$ws.Range("E21").Value = 0.666666666666667
$ws.Range("F21").Value = 0.884615384615385
$ws.Range("G21").Value = 0.102564102564103
$ws.Range("H21").Value = 3.05128205128205
$ws.Range("I21").Value = '{''Levitate'': 99, ''Force'': [59, 14, 26, 19], ''Boost'': [80, 46, 9, 25]}'
$ws.Range("B22").Value = 21
$ws.Range("C22").Value = 0.0476190476190476
$ws.Range("D22").Value = 3.47619047619048
$ws.Range("E22").Value = 0.380952380952381
$ws.Range("F22").Value = 0.428571428571429
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 2.38095238095238
$ws.Range("I22").Value = '{''Levitate'': 24, ''Force'': [9, 4, 4, 1], ''Boost'': [17, 6, 7, 4]}'
$ws.Range("B23").Value = 78
$ws.Range("C23").Value = 0.0128205128205128
$ws.Range("D23").Value = 5.01282051282051
$ws.Range("E23").Value = 0.807692307692308
$ws.Range("F23").Value = 0.769230769230769
$ws.Range("G23").Value = 0.0641025641025641
$ws.Range("H23").Value = 2.66666666666667
$ws.Range("I23").Value = '{''Levitate'': 97, ''Force'': [38, 10, 17, 11], ''Boost'': [73, 35, 17, 21]}'

Write-Output "Prop bets columns resequenced"
